# Fix issue with 0201 components in the BOM file
# (0603 were chosen by jlcpcb auto tool by default).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: C12,C14,... caps -> 0201X104M100CT / 0201 / C2443222
$ws.Range("A3").Value = "0201X104M100CT"
$ws.Range("D3").Value = "0201"
$ws.Range("E3").Value = "C2443222"

# Row 4: C5,C6 caps -> RF03N150F250CT / 0201 / C3896518
$ws.Range("A4").Value = "RF03N150F250CT"
$ws.Range("D4").Value = "0201"
$ws.Range("E4").Value = "C3896518"

# Row 7: R4 -> WR02X102JAL / 0201 / C170199
$ws.Range("A7").Value = "WR02X102JAL"
$ws.Range("D7").Value = "0201"
$ws.Range("E7").Value = "C170199"

# Row 8: R5,R6,R7 -> RM02F5101CT / 0201 / C4153144
$ws.Range("A8").Value = "RM02F5101CT"
$ws.Range("D8").Value = "0201"
$ws.Range("E8").Value = "C4153144"

# Row 9: R2,R3 -> RTT0127R0FTH / 0201 / C158649
$ws.Range("A9").Value = "RTT0127R0FTH"
$ws.Range("D9").Value = "0201"
$ws.Range("E9").Value = "C158649"

# Reflect the last active cell selection seen in the saved file
$ws.Range("D6").Select()
